$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to Text format so numeric-looking strings
# (e.g. values with trailing zeros, multi-dot thousand separators, or tiny
# decimals) are preserved exactly as text, matching the source data feed.
$ws.Range("D2:D51").NumberFormat = "@"

# Update Price (D) and Volume(1h) (E) columns with the refreshed market data
$ws.Range("D2").Value = "20.188.68"
$ws.Range("E2").Value = "  +2.33%  "
$ws.Range("D3").Value = "1.433.08"
$ws.Range("E3").Value = "  +3.43%  "
$ws.Range("D4").Value = "1.009"
$ws.Range("E4").Value = "  +0.57%  "
$ws.Range("D5").Value = "0.9094"
$ws.Range("E5").Value = "  -9.34%  "
$ws.Range("D6").Value = "276.88"
$ws.Range("E6").Value = "  +3.69%  "
$ws.Range("D7").Value = "0.3640"
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").Value = "0.3092"
$ws.Range("E8").Value = "  +2.02%  "
$ws.Range("D9").Value = "38.94"
$ws.Range("E9").Value = "  +3.28%  "
$ws.Range("D10").Value = "1.016"
$ws.Range("E10").Value = "  +4.59%  "
$ws.Range("D11").Value = "0.06501"
$ws.Range("E11").Value = "  +1.64%  "
$ws.Range("D12").Value = "1.004"
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("D13").Value = "5.345"
$ws.Range("E13").Value = "  +1.46%  "
$ws.Range("D14").Value = "17.41"
$ws.Range("E14").Value = "  +5.88%  "
$ws.Range("D15").Value = "6.027"
$ws.Range("E15").Value = "  -0.19%  "
$ws.Range("D16").Value = "0.00001011"
$ws.Range("E16").Value = "  +2.47%  "
$ws.Range("D17").Value = "1.439.29"
$ws.Range("E17").Value = "  +3.67%  "
$ws.Range("D18").Value = "0.9428"
$ws.Range("E18").Value = "  -6.05%  "
$ws.Range("D19").Value = "0.05632"
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("D20").Value = "67.29"
$ws.Range("E20").Value = "  -3.43%  "
$ws.Range("D21").Value = "5.357"
$ws.Range("E21").Value = "  -2.46%  "
$ws.Range("D22").Value = "14.29"
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").Value = "10.75"
$ws.Range("E23").Value = "  +2.42%  "
$ws.Range("D24").Value = "2.224"
$ws.Range("E24").Value = "  -1.09%  "
$ws.Range("D25").Value = "20.261.10"
$ws.Range("E25").Value = "  +2.67%  "
$ws.Range("D26").Value = "2.142"
$ws.Range("E26").Value = "  -0.67%  "
$ws.Range("D27").Value = "136.82"
$ws.Range("E27").Value = "  +0.62%  "
$ws.Range("D28").Value = "16.88"
$ws.Range("E28").Value = "  +2.44%  "
$ws.Range("D29").Value = "1.591.99"
$ws.Range("E29").Value = "  +3.10%  "
$ws.Range("D30").Value = "109.43"
$ws.Range("E30").Value = "  +2.02%  "
$ws.Range("D31").Value = "3.882"
$ws.Range("E31").Value = "  +1.47%  "
$ws.Range("D32").Value = "0.7957"
$ws.Range("E32").Value = "  +0.85%  "
$ws.Range("D33").Value = "4.746"
$ws.Range("E33").Value = "  -8.69%  "
$ws.Range("D34").Value = "0.07646"
$ws.Range("E34").Value = "  +1.12%  "
$ws.Range("D35").Value = "0.05913"
$ws.Range("E35").Value = "  +5.87%  "
$ws.Range("D36").Value = "1.438"
$ws.Range("E36").Value = "  +12.30%  "
$ws.Range("D37").Value = "1.131"
$ws.Range("E37").Value = "  +8.32%  "
$ws.Range("D38").Value = "4.610"
$ws.Range("E38").Value = "  -1.00%  "
$ws.Range("D39").Value = "0.01976"
$ws.Range("E39").Value = "  -1.73%  "
$ws.Range("D40").Value = "10.11"
$ws.Range("E40").Value = "  +1.99%  "
$ws.Range("D41").Value = "0.1826"
$ws.Range("E41").Value = "  -1.97%  "
$ws.Range("D42").Value = "0.9204"
$ws.Range("E42").Value = "  -8.25%  "
$ws.Range("D43").Value = "7.023"
$ws.Range("E43").Value = "  -14.03%  "
$ws.Range("D46").Value = "11.98"
$ws.Range("E46").Value = "  +1.57%  "
$ws.Range("D47").Value = "118.04"
$ws.Range("E47").Value = "  +9.15%  "
$ws.Range("D48").Value = "0.5086"
$ws.Range("E48").Value = "  +2.70%  "
$ws.Range("D49").Value = "1.745"
$ws.Range("E49").Value = "  +1.78%  "
$ws.Range("D50").Value = "0.06301"
$ws.Range("E50").Value = "  +4.33%  "
$ws.Range("D51").Value = "0.9880"
$ws.Range("E51").Value = "  -1.42%  "

# PancakeSwap overtook TheSandbox in the rankings -> rows 44/45 swap places
$ws.Range("B44").Value = "PancakeSwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D44").Value = "3.505"
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "0.5198"
$ws.Range("E45").Value = "  +0.89%  "
